$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116; Excel shifts rows 116:188 down to 117:189
# and the new row inherits the formatting (incl. the date number format on
# column D) from the row above, just like a native Excel "Insert Row" does.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new weekly record.
$ws.Range("A116").Value = 4
$ws.Range("B116").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C116").Value = "Los Lagos"
$ws.Range("D116").Value = 44582
$ws.Range("E116").Value = 10
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100103
$ws.Range("H116").Value = "Frutos de hueso (carozo)"
$ws.Range("I116").Value = 100103004
$ws.Range("J116").Value = "Durazno"
$ws.Range("K116").Value = "Carson"
$ws.Range("L116").Value = "Primera"
$ws.Range("M116").Value = 500
$ws.Range("N116").Value = 17000
$ws.Range("O116").Value = 18000
$ws.Range("P116").Value = 17500
$ws.Range("Q116").Value = '$/caja 15 kilos empedrada'
$ws.Range("R116").Value = "Región de O'Higgins"
$ws.Range("S116").Value = 1167
$ws.Range("T116").Value = 15
